# Update the date line (2024-05-08 Wednesday -> 2024-05-09 Thursday)
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-05-08 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-05-09 Thursday", 2)

# Update each multiplication problem cell in the first (and only) table.
# Addressing cells directly (row/column) avoids any ambiguity that could
# arise from a plain text find/replace, since some "new" values equal
# "old" values used elsewhere in the table (e.g. 50x43=).
$tbl = $d.Tables.Item(1)

$newValues = @{
    "1,1"  = "35×36=";  "1,2"  = "75×96=";  "1,3"  = "81×17=";  "1,4"  = "17×67=";  "1,5"  = "22×64=";
    "5,1"  = "15×52=";  "5,2"  = "17×29=";  "5,3"  = "75×33=";  "5,4"  = "51×73=";  "5,5"  = "19×95=";
    "10,1" = "89×27=";  "10,2" = "50×43=";  "10,3" = "33×13=";  "10,4" = "13×96=";  "10,5" = "73×74=";
    "15,1" = "23×87=";  "15,2" = "21×73=";  "15,3" = "97×78=";  "15,4" = "11×90=";  "15,5" = "52×79=";
    "20,1" = "22×48=";  "20,2" = "30×72=";  "20,3" = "56×48=";  "20,4" = "54×80=";  "20,5" = "94×15=";
}

$tableRows = @(1, 5, 10, 15, 20)
foreach ($tblRow in $tableRows) {
    for ($colIdx = 1; $colIdx -le 5; $colIdx++) {
        $key = "$tblRow,$colIdx"
        $cell = $tbl.Cell($tblRow, $colIdx)
        $cell.Range.Text = $newValues[$key]
    }
}
